# Generate Report for Handoff
# Moves the two tracked files from their previous handback state to a
# fresh "Ready for handoff" state: new source/target GUID-named files,
# new handoff timestamps, and clears the (now stale) target/handback
# columns since nothing has been handed back yet for the new cycle.

$wb = $excel.ActiveWorkbook

$oldMd1 = "48c498d6-18be-4962-98d1-413d890aaeb6.md"
$oldMd2 = "f6315fb3-99db-4331-9f89-825288b7833e.md"
$newMd1 = "cd3d2eed-9657-46d0-a4ae-8a7a25503f11.md"
$newMd2 = "ffffc78e3aa7-9863-43d8-afdb-b332caab3b73.md"

$newStatus = "Ready for handoff"

$newHandoffDate = "2016-03-23 09:14:08"

$newTargetZh = "cd3d2eed-9657-46d0-a4ae-8a7a25503f11.72c9c9f2dede7c896e651960276eab8bc9184b97.zh-cn.xlf"
$newTargetDe = "cd3d2eed-9657-46d0-a4ae-8a7a25503f11.72c9c9f2dede7c896e651960276eab8bc9184b97.de-de.xlf"

$newHandoffDatetimeZh = "2016-03-23 09:14:04"
$newHandbackDatetime = "0001-01-01 00:00:00"

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Item(2,1).Value2 = $newMd1
$ov.Cells.Item(2,2).Value2 = $newStatus
$ov.Cells.Item(2,3).Value2 = $newStatus
$ov.Cells.Item(2,4).Value2 = $newHandoffDate
$ov.Cells.Item(3,1).Value2 = $newMd2
$ov.Cells.Item(3,2).Value2 = $newStatus
$ov.Cells.Item(3,3).Value2 = $newStatus
$ov.Cells.Item(3,4).Value2 = $newHandoffDate

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address(0,0)
    if ($addr -eq "A2") {
        $hl.TextToDisplay = $newMd1
    } elseif ($addr -eq "A3") {
        $hl.TextToDisplay = $newMd2
    }
}

# ---- zh-cn and de-de sheets share the same shape ----
function Update-LangSheet($ws, $newTarget, $newHandoffDatetime) {
    $ws.Cells.Item(2,1).Value2 = $newMd1
    $ws.Cells.Item(2,3).Value2 = $newStatus
    $ws.Cells.Item(2,4).Value2 = $newTarget
    $ws.Cells.Item(2,5).Value2 = $newHandoffDatetime
    $ws.Cells.Item(2,6).Clear()
    $ws.Cells.Item(2,7).Clear()
    $ws.Cells.Item(2,8).Value2 = $newHandbackDatetime

    $ws.Cells.Item(3,1).Value2 = $newMd2
    $ws.Cells.Item(3,3).Value2 = $newStatus
    $ws.Cells.Item(3,4).Value2 = $newTarget
    $ws.Cells.Item(3,5).Value2 = $newHandoffDatetime
    $ws.Cells.Item(3,6).Clear()
    $ws.Cells.Item(3,7).Clear()
    $ws.Cells.Item(3,8).Value2 = $newHandbackDatetime

    # Drop the hyperlinks anchored on the now-empty Latest Target File /
    # Latest Handback File cells (F/G). Re-resolve the collection fresh
    # for each single delete -- deleting while a stale reference to a
    # sibling hyperlink is held can silently no-op on this host.
    $deadAnchors = @("F2", "G2", "F3", "G3")
    foreach ($anchor in $deadAnchors) {
        $victim = $null
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address(0,0) -eq $anchor) {
                $victim = $hl
            }
        }
        if ($victim -ne $null) {
            $victim.Delete()
        }
    }

    # Refresh display text + target on the hyperlinks that remain.
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address(0,0)
        if ($addr -eq "A2") {
            $hl.TextToDisplay = $newMd1
        } elseif ($addr -eq "A3") {
            $hl.TextToDisplay = $newMd2
        } elseif ($addr -eq "D2" -or $addr -eq "D3") {
            $hl.TextToDisplay = $newTarget
        }
    }
}

$zh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zh $newTargetZh $newHandoffDatetimeZh

$de = $wb.Worksheets.Item("de-de")
Update-LangSheet $de $newTargetDe $newHandoffDate
